$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its original text formatting so numeric-looking
# strings like "1.00" or "599.16" are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.582.17"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "3.567.70"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "599.16"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "140.33"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("D7").Value = "3.567.64"
$ws.Range("E7").Value = "  +3.39%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").Value = "7.14"
$ws.Range("E11").Value = "  -5.29%  "
$ws.Range("D12").Value = "0.393"
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("D13").Value = "4.174.45"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "3.568.18"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "27.08"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "65.471.43"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "10.30"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("D20").Value = "5.86"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").Value = "14.27"
$ws.Range("E21").Value = "  +3.62%  "
$ws.Range("D22").Value = "396.44"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "0.571"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("D24").Value = "3.714.04"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").Value = "74.73"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +10.44%  "
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +7.12%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").Value = "3.586.20"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").Value = "23.96"
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("D36").Value = "1.26"
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("D37").Value = "7.06"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D39").Value = "168.39"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "5.00"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("D41").Value = "0.0803"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "26.71"
$ws.Range("E43").Value = "  +15.57%  "
$ws.Range("D44").Value = "43.02"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "4.44"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "1.71"
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").Value = "2.438.90"
$ws.Range("E49").Value = "  +10.48%  "
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  +4.53%  "
$ws.Range("E51").Value = "  +1.97%  "

# Row 34/35 coin entries swapped (Kaspa <-> USDe) with updated price/volume data
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  +0.55%  "

# Restore default style on column D so no stray formatting is left behind
$ws.Range("D2:D51").Style = "Normal"
